# newly added iAuthor TC's
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-9 (8 rows), columns A-H
$data = @(
    @("CIO12057", 20123201, "candidatepfm12057", "Elumina@179", "MR", "Candidate", "PFTMM", "Candidate"),
    @("CIO12056", 20123200, "candidatepfm12056", "Elumina@178", "MR", "Candidate", "PFTMM", "Candidate"),
    @("CIO12055", 20123199, "candidatepfm12055", "Elumina@177", "MR", "Candidate", "PFTMM", "Candidate"),
    @("CIO12054", 20123198, "candidatepfm12054", "Elumina@176", "MR", "Candidate", "PFTMM", "Candidate"),
    @("CIO12053", 20123197, "candidatepfm12053", "Elumina@175", "MR", "Candidate", "PFTMM", "Candidate"),
    @("CIO12052", 20123196, "candidatepfm12052", "Elumina@174", "MR", "Candidate", "PFTMM", "Candidate"),
    @("CIO12051", 20123195, "candidatepfm12051", "Elumina@173", "MR", "Candidate", "PFTMM", "Candidate"),
    @("CIO12050", 20123194, "candidatepfm12050", "Elumina@172", "MR", "Candidate", "PFTMM", "Candidate")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}

# Match the thin-border look of the existing data rows for the newly added ones (rows 7-9)
$ws.Range("A7:H9").Borders.LineStyle = 1

# Extend the selection to cover the new full range, mirroring the expanded table
$null = $ws.Range("A1:H9").Select()
